# Edit script for national_defense_pe_original_data.xlsx
# Commit: "Update trade day data."
#
# Summary of changes:
# 1. Row 4 (HIGH/LOW/WIDTH/AVE/STD summary row): AVE (E4) and STD (F4) are
#    recalculated because two new trading-day observations were added to the
#    "pe" series (replacing two duplicate placeholder rows at the bottom).
# 2. Rows 6-248 (data rows): the constant "ave_pe" (D), "plus_1sigma" (E) and
#    "minus_1sigma" (F) columns are refreshed to reflect the new AVE/STD.
# 3. Rows 249-250: previously duplicated the data from rows 247-248 (same
#    date serials 44434/44435). They are now replaced with real data for the
#    next two trading days (44438 / 44439, i.e. 2021-08-30 / 2021-08-31),
#    including the date, pe value, the refreshed ave_pe/plus_1sigma/minus_1sigma
#    and all per-security weight columns (G:CA).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Row 4 summary stats: AVE (E4) and STD (F4)
# ---------------------------------------------------------------------------
$ws.Range("E4").Value2 = [double]79.29551896260523
$ws.Range("F4").Value2 = [double]61.07447281236171

# ---------------------------------------------------------------------------
# 2) Rows 6-248: refresh ave_pe (D), plus_1sigma (E), minus_1sigma (F)
#    These are identical constants repeated down the whole data range, so we
#    build one small 2D array and write it across the full D6:F248 block.
# ---------------------------------------------------------------------------
$firstDataRow = 6
$lastOldDataRow = 248
$nRows = $lastOldDataRow - $firstDataRow + 1

$dEFArray = New-Object 'object[,]' $nRows,3
for ($r = 0; $r -lt $nRows; $r++) {
    $dEFArray[$r,0] = [double]70.18499588748347
    $dEFArray[$r,1] = [double]79.29551896260523
    $dEFArray[$r,2] = [double]61.07447281236171
}
$ws.Range("D" + $firstDataRow + ":F" + $lastOldDataRow).Value2 = $dEFArray

# ---------------------------------------------------------------------------
# 3) Rows 249 & 250: replace duplicated placeholder rows with the real data
#    for the next two trading days across every column B:CA.
# ---------------------------------------------------------------------------
$csv249 = "44438,73.875943646433,70.18499588748347,79.29551896260523,61.07447281236171,0.414898585588074,1.143073714287487,0.1843338305348915,0.06259314114835723,0.0972601506876168,0.2278806517479913,1.432565607013954,3.338065539440373,1.260059090833572,0.7633766172068638,0.8281787861725441,0.4084457542596633,0.08707257311822542,0.1726524946694332,0.4186177651233256,0.9079912205076323,0.9131308099567517,0.5100282818613954,0.162771016144381,0.9807820529238803,0.4672967218207355,0.6482873332430446,0.1390847225025745,2.535387299372488,-0.03027202245153993,0.1462524224826671,0.7514925827394839,0.6437730052841714,0.8309450522758471,6.500645196820455,-0.2444967356387809,0.4995386754414862,0.220599653526883,1.703255630802023,0.2325080526864387,0.8660747888972462,0.3296781953487565,1.35557565731809,0.2089537328645305,0.8836812464877435,-12.34423367072106,0.2575040246549273,0.3920658695025171,3.148717957024914,1.186218636922694,0.06863266489079513,0.8185475932010872,1.092008549020755,0.2216953407631566,0.1534489374338729,0.256528977637373,0.5949403655855,0.3520291190484699,0.9089956874860022,0.139438875696695,0.7921157532410871,3.104117236916785,1.476743972448429,1.798811981344264,0.9894520686651979,4.561721690648091,-0.03889445870949414,-0.002552373801944852,1.654737382112033,0.2979618640726581,10.05191479425953,13.39502889601698,0.3750866529981813,2.613711408923936,0.1155184185206094,0.5950227645388702,1.022630450557278,0.8242373204820674"
$parts249 = $csv249 -split ','
$n249 = $parts249.Length
$arr249 = New-Object 'object[,]' 1,$n249
for ($i = 0; $i -lt $n249; $i++) {
    $arr249[0,$i] = [double]$parts249[$i]
}
$ws.Range("B249:CA249").Value2 = $arr249

$csv250 = "44439,70.6616020481805,70.18499588748347,79.29551896260523,61.07447281236171,0.4051724541764333,1.116501078669658,0.1841223415138565,0.0609427753602363,0.1008222262728551,0.2282949802027998,1.459279011867268,3.495333391650508,1.238745823745302,0.7554413093533526,0.8319902600117006,0.405633945673752,0.08642459116350155,0.1697436425259534,0.4120768625439651,0.9271453372047204,0.9068751070937248,0.4935073394817977,0.1655687991045788,0.9601436960003451,0.4666180047165208,0.6291271043127384,0.1333561728195923,2.543950627189562,-0.032219640095883,0.1531819222241976,0.7540156242298575,0.6258508210374573,0.8230111402488092,6.364170640450669,-0.2333612407523648,0.4986474378494897,0.2188427623341878,1.703255630802023,0.2315015242862359,0.8362101409772592,0.3540168540922775,1.337096534407987,0.2028536007188775,0.9002514467127319,-12.96738008215186,0.2495740461310947,0.3844594179899234,3.120084760064713,1.173121994114706,0.06795581416408626,0.8111160584574161,1.140053667960742,0.2201753923185344,0.1527464001406401,0.25631784270086,0.6062295565840842,0.3459435092905877,0.9103977374047522,0.1378016286283371,0.7879705233374256,3.244665882834134,1.479390467134738,1.914165321386028,1.001387093571634,4.470043653398048,-0.03876872230062903,-0.002677490163966699,1.39754234552705,0.2974857346200783,7.938047615719447,13.322623334493,0.3731946647501777,2.615832354143098,0.1124001299264448,0.5355779227310761,0.8836160206990852,0.8382795671332804"
$parts250 = $csv250 -split ','
$n250 = $parts250.Length
$arr250 = New-Object 'object[,]' 1,$n250
for ($i = 0; $i -lt $n250; $i++) {
    $arr250[0,$i] = [double]$parts250[$i]
}
$ws.Range("B250:CA250").Value2 = $arr250
